$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.412.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.507.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.44%  "
$ws.Range("E7").Value = "  +1.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.55%  "
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.898.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.502.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.351.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.65%  "
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "249.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.04%  "
$ws.Range("E31").Value = "  +9.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0802"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.92%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  +6.44%  "
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.24%  "
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.007.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E46").Value = "  +5.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
